# UPSMS-Scorecard.xlsx -- "Minor changes + Uploaded Deliverables"
#
# What changed (per the target diff):
#   - User sheet:       F4 and F5 ("Complied?" column) now say "Yes";
#                        selection moves to F5; sheet becomes the active tab.
#   - Admin sheet:       F4 and F5 ("Complied?" column) now say "Yes";
#                        selection moves to F5.
#   - Signatory sheet:   stops being the active tab (selection/scroll untouched).
#   - Workbook:          active tab becomes the User sheet (tab 1 / index 0).

$wb = $excel.ActiveWorkbook

$wsUser      = $wb.Worksheets.Item("User")
$wsSignatory = $wb.Worksheets.Item("Signatory")
$wsAdmin     = $wb.Worksheets.Item("Admin")

# --- User sheet: mark rows 4 & 5 of the "Complied?" column (F) as "Yes" ---
$wsUser.Cells.Item(4, 6).Value = "Yes"
$wsUser.Cells.Item(5, 6).Value = "Yes"

# --- Admin sheet: same "Uploaded Deliverables" mark for rows 4 & 5 ---
$wsAdmin.Cells.Item(4, 6).Value = "Yes"
$wsAdmin.Cells.Item(5, 6).Value = "Yes"

# Admin's selection ends on F5
$wsAdmin.Range("F5").Select()

# Signatory sheet is left as-is (its selection/scroll position don't change);
# simply activating User afterwards below moves "tabSelected" off of it.

# User sheet ends up selected at F5 and becomes the active/visible tab
$wsUser.Range("F5").Select()
$wsUser.Activate()
